$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Capture the existing detailed header row (row 2) before it gets overwritten ---
$row2Vals = @{}
for ($c = 2; $c -le 23; $c++) {
    $row2Vals[$c] = $ws.Cells.Item(2, $c).Value2
}

# --- 2. Unmerge the grouped header cells in row 1 ---
$ws.Range("H1:L1").UnMerge()
$ws.Range("M1:P1").UnMerge()
$ws.Range("Q1:S1").UnMerge()

# --- 3. Rewrite row 1 with the real column labels (cleaned headers) ---
$ws.Cells.Item(1, 1).Value = "Player ID"
for ($c = 2; $c -le 23; $c++) {
    $ws.Cells.Item(1, $c).Value = $row2Vals[$c]
}
$ws.Cells.Item(1, 7).Value = "90s"
$ws.Cells.Item(1, 13).Value = "Cha"

# --- 4. Restore row 2 to its original (pre-edit) detailed header content, now as a hidden duplicate ---
for ($c = 2; $c -le 23; $c++) {
    $ws.Cells.Item(2, $c).Value = $row2Vals[$c]
}
$ws.Rows(2).Hidden = $true

# --- 5. Insert a new, blank hidden row 3 (no shifting of existing rows) ---
$ws.Rows(3).Hidden = $true

# --- 6. Fill in the missing "Tkl%" (O) values with 0 for the rows that lacked them ---
$oRows = @(4,5,6,8,9,10,11,12,13,18,19)
foreach ($r in $oRows) {
    $ws.Cells.Item($r, 15).Value = 0
}

# --- 7. Hide the totals row ---
$ws.Rows(20).Hidden = $true

# --- 8. Leave the selection where the author left it ---
$ws.Range("O21").Select() | Out-Null
